# Inserts a new weekly price record row at row 455 of the "Zanahoria" sheet.
# This pushes the previous rows 455-571 down to 456-572, and the new row 455
# gets fresh data for the "Macroferia Regional de Talca" / Zanahoria series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 455 (existing rows 455.. shift down by one)
$ws.Rows.Item(455).Insert()

# Populate the newly inserted row with the new record
$ws.Cells.Item(455, 1).Value = 5
$ws.Cells.Item(455, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(455, 3).Value = "Maule"
$ws.Cells.Item(455, 4).Value = 45135
$ws.Cells.Item(455, 5).Value = 7
$ws.Cells.Item(455, 6).Value = 100114013
$ws.Cells.Item(455, 7).Value = "Zanahoria"
$ws.Cells.Item(455, 8).Value = "Sin especificar"
$ws.Cells.Item(455, 9).Value = "Primera"
$ws.Cells.Item(455, 10).Value = 500
$ws.Cells.Item(455, 11).Value = 5000
$ws.Cells.Item(455, 12).Value = 5000
$ws.Cells.Item(455, 13).Value = 5000
$ws.Cells.Item(455, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(455, 15).Value = "Región de Ñuble"
$ws.Cells.Item(455, 16).Value = 250
$ws.Cells.Item(455, 17).Value = 20
$ws.Cells.Item(455, 18).Value = "Hortaliza"
